$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace placeholder hotel names with real hotel names
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("B3").Value = "Grand Plaza Apartments"
